$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 405, shifting existing rows 405..481 down to 407..483.
$ws.Range("A405:A406").EntireRow.Insert()

# Populate new row 405 with its final values.
$ws.Range("A405").Value = 4
$ws.Range("B405").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C405").Value = "Los Lagos"
$ws.Range("D405").Value = 44637
$ws.Range("E405").Value = 10
$ws.Range("F405").Value = 100112004
$ws.Range("G405").Value = "Cebolla"
$ws.Range("H405").Value = "Morada(o)"
$ws.Range("I405").Value = "1a (cosecha)"
$ws.Range("J405").Value = 120
$ws.Range("K405").Value = 14000
$ws.Range("L405").Value = 14000
$ws.Range("M405").Value = 14000
$ws.Range("N405").Value = "$/malla 18 kilos"
$ws.Range("O405").Value = "Región de O'Higgins"
$ws.Range("P405").Value = 778
$ws.Range("Q405").Value = 18
$ws.Range("R405").Value = "Hortaliza"

# Populate new row 406 with its final values.
$ws.Range("A406").Value = 4
$ws.Range("B406").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C406").Value = "Los Lagos"
$ws.Range("D406").Value = 44637
$ws.Range("E406").Value = 10
$ws.Range("F406").Value = 100112004
$ws.Range("G406").Value = "Cebolla"
$ws.Range("H406").Value = "Sin especificar"
$ws.Range("I406").Value = "1a (cosecha)"
$ws.Range("J406").Value = 400
$ws.Range("K406").Value = 8000
$ws.Range("L406").Value = 8500
$ws.Range("M406").Value = 8250
$ws.Range("N406").Value = "$/malla 18 kilos"
$ws.Range("O406").Value = "Región de O'Higgins"
$ws.Range("P406").Value = 458
$ws.Range("Q406").Value = 18
$ws.Range("R406").Value = "Hortaliza"

# Ensure the date cells keep the expected date number format used throughout column D.
$ws.Range("D405:D406").NumberFormat = "YYYY-MM-DD HH:MM:SS"
